# SQA signoff for conflicts (FileIndex.xlsx)
#
# 1. Row 16's description is updated: the old "rejected - duplicate classes"
#    message is replaced by the new "accepted - duplicate classes, line
#    three overwrites line one" message (the old shared string disappears
#    from the sharedStrings table since nothing else references it, and the
#    new text is appended as a new shared string - Excel handles this
#    automatically when the cell value is set).
# 2. Row 23 is taller now (more text wraps), so its custom row height grows
#    from 83.25 to 108.75.
# 3. The sheet view had scrolled to show row 19 at the top with B23 selected;
#    now it shows row 20 at the top with A21 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Input is accepted. Line one and line three have duplicate classes so line three will overwrite line one in the database."

$ws.Rows.Item(23).RowHeight = 108.75

$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A21").Select()
